# Thrombophilia.xlsx — "Refined metadata to be additional tab"
#
# 1. Update the panel_query_time / time_taken timestamps on the "data" sheet
#    (re-run recorded at 2021-10-05 14:22:57.xxxxxx instead of 13:42:27.xxxxxx).
# 2. Add a new "metadata" worksheet (after "data") carrying the panel-level
#    metadata that used to only implicitly exist, with the same header/
#    formatting style as the "data" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. Refresh the per-row "time_taken" timestamps on the data sheet ------
$ws.Range("F2").Value  = "2021-10-05 14:22:57.868717"
$ws.Range("F3").Value  = "2021-10-05 14:22:57.868725"
$ws.Range("F4").Value  = "2021-10-05 14:22:57.868728"
$ws.Range("F5").Value  = "2021-10-05 14:22:57.868731"
$ws.Range("F6").Value  = "2021-10-05 14:22:57.868734"
$ws.Range("F7").Value  = "2021-10-05 14:22:57.868737"
$ws.Range("F8").Value  = "2021-10-05 14:22:57.868739"
$ws.Range("F9").Value  = "2021-10-05 14:22:57.868742"
$ws.Range("F10").Value = "2021-10-05 14:22:57.868745"
$ws.Range("F11").Value = "2021-10-05 14:22:57.868747"
$ws.Range("F12").Value = "2021-10-05 14:22:57.868750"
$ws.Range("F13").Value = "2021-10-05 14:22:57.868752"
$ws.Range("F14").Value = "2021-10-05 14:22:57.868755"
$ws.Range("F15").Value = "2021-10-05 14:22:57.868757"
$ws.Range("F16").Value = "2021-10-05 14:22:57.868760"
$ws.Range("F17").Value = "2021-10-05 14:22:57.868762"
$ws.Range("F18").Value = "2021-10-05 14:22:57.868765"
$ws.Range("F19").Value = "2021-10-05 14:22:57.868768"
$ws.Range("F20").Value = "2021-10-05 14:22:57.868770"
$ws.Range("F21").Value = "2021-10-05 14:22:57.868773"
$ws.Range("F22").Value = "2021-10-05 14:22:57.868775"

# --- 2. Add the new "metadata" sheet, placed right after "data" -----------
$metaSheet = $wb.Worksheets.Add($null, $ws)
$metaSheet.Name = "metadata"

# Reuse the same header style (bold / bordered / centered / top-aligned)
# that the "data" sheet's header row + first column already carry, instead
# of inventing a brand-new style.
$ws.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$ws.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Thrombophilia"
$metaSheet.Range("C2").Value = 516
# "1.20" must stay textual (not become the number 1.2) — enter it with a
# leading apostrophe to force text, then re-paste a plain (unstyled) format
# over it so the cell doesn't keep the quote-prefix style.
$metaSheet.Range("D2").Value = "'1.20"
$ws.Range("B2").Copy()
$metaSheet.Range("D2").PasteSpecial(-4122)
$metaSheet.Range("E2").Value = "2021-07-27T10:21:19.668610Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:22:57.865053"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/516/?format=json"

$excel.CutCopyMode = $false

# Keep "data" as the active sheet, matching the original workbook's state.
$ws.Activate()

Write-Host "done"
